$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $r = $ws.Range($cell)
    $r.Value = "'" + $text
    $r.ClearFormats()
}

$ws.Range('D2').Value = '68.349.07'
$ws.Range('E2').Value = '  +2.31%  '
$ws.Range('D3').Value = '3.134.92'
$ws.Range('E3').Value = '  +1.92%  '
$ws.Range('E4').Value = '  +0.09%  '
Set-TextValue 'D5' '577.59'
$ws.Range('E5').Value = '  +0.21%  '
Set-TextValue 'D6' '180.98'
$ws.Range('E6').Value = '  +5.95%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '3.133.58'
$ws.Range('E8').Value = '  +1.93%  '
$ws.Range('E9').Value = '  +1.61%  '
Set-TextValue 'D10' '6.55'
$ws.Range('E10').Value = '  +2.74%  '
Set-TextValue 'D11' '0.153'
$ws.Range('E11').Value = '  +1.91%  '
Set-TextValue 'D12' '0.471'
$ws.Range('E12').Value = '  +0.87%  '
Set-TextValue 'D13' '0.0000242'
$ws.Range('E13').Value = '  +1.16%  '
Set-TextValue 'D14' '36.84'
$ws.Range('E14').Value = '  +2.85%  '
$ws.Range('B15').Value = 'TRON'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D15' '0.122'
$ws.Range('E15').Value = '  +1.00%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '68.332.47'
$ws.Range('E16').Value = '  +2.32%  '
$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = '3.663.94'
$ws.Range('E17').Value = '  +2.12%  '
Set-TextValue 'D18' '7.12'
$ws.Range('E18').Value = '  +1.73%  '
$ws.Range('D19').Value = '3.143.03'
$ws.Range('E19').Value = '  +2.13%  '
Set-TextValue 'D20' '16.68'
$ws.Range('E20').Value = '  -1.87%  '
Set-TextValue 'D21' '488.88'
$ws.Range('E21').Value = '  -0.72%  '
Set-TextValue 'D22' '0.698'
$ws.Range('E22').Value = '  +1.34%  '
Set-TextValue 'D23' '7.80'
$ws.Range('E23').Value = '  +1.10%  '
Set-TextValue 'D24' '83.94'
$ws.Range('E24').Value = '  +1.32%  '
Set-TextValue 'D25' '13.00'
$ws.Range('E25').Value = '  +2.47%  '
Set-TextValue 'D26' '2.32'
$ws.Range('E26').Value = '  +5.50%  '
Set-TextValue 'D27' '10.58'
$ws.Range('E27').Value = '  +4.19%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('E29').Value = '  +4.04%  '
Set-TextValue 'D30' '2.36'
$ws.Range('E30').Value = '  +3.92%  '
$ws.Range('E31').Value = '  +1.08%  '
Set-TextValue 'D32' '28.38'
$ws.Range('E32').Value = '  +2.74%  '
Set-TextValue 'D33' '0.113'
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('D34').Value = '0.0₃0956'
$ws.Range('E34').Value = '  +4.43%  '
Set-TextValue 'D35' '1.00'
$ws.Range('E35').Value = '  +0.15%  '
Set-TextValue 'D36' '49.00'
$ws.Range('E36').Value = '  +3.86%  '
Set-TextValue 'D37' '5.65'
$ws.Range('E37').Value = '  +1.01%  '
Set-TextValue 'D38' '0.957'
$ws.Range('E38').Value = '  +0.72%  '
Set-TextValue 'D39' '0.322'
$ws.Range('E39').Value = '  +7.34%  '
Set-TextValue 'D40' '2.06'
$ws.Range('E40').Value = '  +4.74%  '
Set-TextValue 'D41' '49.13'
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('E42').Value = '  +1.18%  '
Set-TextValue 'D43' '8.43'
$ws.Range('E43').Value = '  +1.34%  '
Set-TextValue 'D44' '2.70'
$ws.Range('E44').Value = '  +8.28%  '
Set-TextValue 'D45' '395.77'
$ws.Range('E45').Value = '  +7.62%  '
$ws.Range('D46').Value = '2.785.22'
$ws.Range('E46').Value = '  +1.08%  '
$ws.Range('E47').Value = '  +10.22%  '
$ws.Range('E48').Value = '  +1.11%  '
Set-TextValue 'D49' '135.63'
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('E51').Value = '  +8.66%  '
